$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.802.06'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.55%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.368.11'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.78%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.44'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.86%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.55'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.09%  '

$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.365.55'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.83%  '

$ws.Range("E9").Value = '  -1.22%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.63'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.74%  '

$ws.Range("E11").Value = '  -2.58%  '

$ws.Range("E12").Value = '  -1.77%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.942.62'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.78%  '

$ws.Range("E14").Value = '  +0.55%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.82'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.65%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.367.50'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.78%  '

$ws.Range("E17").Value = '  -3.48%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '60.968.64'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.33%  '

$ws.Range("E19").Value = '  -2.00%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.80'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.82%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.40'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.02%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '373.71'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.77%  '

$ws.Range("E23").Value = '  -2.79%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.508.82'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.67%  '

$ws.Range("E25").Value = '  -0.11%  '

$ws.Range("E26").Value = '  -2.22%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '70.96'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.46%  '

$ws.Range("E28").Value = '  +11.81%  '

$ws.Range("E29").Value = '  +9.52%  '

$ws.Range("E30").Value = '  -0.02%  '

$ws.Range("E31").Value = '  -2.87%  '

$ws.Range("E32").Value = '  -2.15%  '

$ws.Range("E33").Value = '  -1.04%  '

$ws.Range("E34").Value = '  +0.00%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.66'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.72%  '

$ws.Range("E36").Value = '  -4.50%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.84'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.54%  '

$ws.Range("E38").Value = '  -1.38%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '164.08'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.91%  '

$ws.Range("E40").Value = '  -3.36%  '

$ws.Range("E42").Value = '  -1.15%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '41.50'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.55%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.69'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.06%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.39'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.11%  '

$ws.Range("E46").Value = '  -2.95%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.14'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.92%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.453.90'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.88%  '

$ws.Range("E49").Value = '  -2.17%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.92'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.06%  '

$ws.Range("E51").Value = '  +3.71%  '
